$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 452, shifting existing rows 452:512 down to 453:513
$ws.Rows.Item(452).Insert()

# Populate the newly inserted row 452 with its data (same shape as the
# other data rows, copied from the row that used to occupy 452 but with
# the updated Fecha / Volumen / Origen values from the commit).
$ws.Range("A452").Value = 4
$ws.Range("B452").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C452").Value = "Los Lagos"
$ws.Range("D452").Value = 45142
$ws.Range("E452").Value = 10
$ws.Range("F452").Value = 100112045
$ws.Range("G452").Value = "Zapallo"
$ws.Range("H452").Value = "Paine"
$ws.Range("I452").Value = "1a (guarda)"
$ws.Range("J452").Value = 1200
$ws.Range("K452").Value = 600
$ws.Range("L452").Value = 600
$ws.Range("M452").Value = 600
$ws.Range("N452").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O452").Value = "Región de O'Higgins"
$ws.Range("P452").Value = 600
$ws.Range("Q452").Value = 1
$ws.Range("R452").Value = "Hortaliza"

# Match the date-style formatting used by the other rows in column D
$ws.Range("D452").NumberFormat = $ws.Range("D453").NumberFormat
